# Updated cryptos list (Price / Volume(1h) columns) - commit: "Updated cryptos list on Tue May 16 20:51:03 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format so numeric-looking strings (e.g. "310.83", "1.0000") are preserved verbatim as text,
# matching the original inlineStr cell type rather than being auto-coerced to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.977.25"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.89%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.823.02"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.83"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.64%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4245"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3654"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07221"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8413"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.11%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.825.44"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.643"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07046"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.269"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.47"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008744"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.12%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.83"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.096.87"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.120"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.78"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.049.67"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.977"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.75"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.216"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.15"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.205"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.50"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08703"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.172"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7336"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.900"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.407"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.0000"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.59%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01934"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05205"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.218"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.874"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.49%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5104"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.505"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.53"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.951"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +6.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4723"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.66"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9995"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06314"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.646"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.47%  "
